# AutoCommit_15 марта 2024 г. 9:58:33_SibNout2023
# Grade entry: three students received a score of 5 (rows 5, 9 and 16),
# and the active cell/selection was left on D9 after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 5

$ws.Range("D9").Value = 5

$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 5

# Leave the selection on D9, matching the author's final cursor position.
$ws.Range("D9").Select()
